$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reassign course "Área" (col A) / "modalidad" (col C) between rows 5 and 6,
#     and rename a couple of area labels, per "Semana 3 | dia 1 | Desarrollo de Datamart OEC" ---

# "Desarrollo de Software" area renamed to "Python"
$ws.Range("A7").Value = "Python"
$ws.Range("A8").Value = "Python"

# Row 5 now holds the Power BI 1 / Presencial course that used to live in row 6
$ws.Range("A5").Value = "Power BI"
$ws.Range("B5").Value = "Power BI 1: Entorno y Publicación"
$ws.Range("C5").Value = "Presencial"
$ws.Range("D5").Value = 24
$ws.Range("E5").Value = 99

# Row 6 now holds the SQL course (area renamed from "Business Intelligence" to "SQL")
$ws.Range("A6").Value = "SQL"
$ws.Range("B6").Value = "SQL Aplicado al Análisis de Datos"
$ws.Range("C6").Value = "En vivo"
$ws.Range("D6").Value = 16
$ws.Range("E6").Value = 90

# "Programación" area renamed to "Programación en R"
$ws.Range("A20").Value = "Programación en R"

# --- Column widths: widen Área / Nombre del curso columns, drop best-fit ---
$ws.Columns.Item(1).ColumnWidth = 20.1667
$ws.Columns.Item(2).ColumnWidth = 76.6667

# --- Row heights: with column B much wider, text no longer wraps, so clear
#     the explicit per-row heights back to the sheet default ---
$ws.Range("A1:E22").Rows.AutoFit()

# --- Selection / view: scroll back to top, select B24 ---
$ws.Range("B24").Select()
